# Sync attendance_reports: normalize "Recorded By" (column G) ordering.
# For every data row, the comma-separated list of recorder names/emails in
# column G is reversed in order (e.g. "a, System" -> "System, a").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Find the "Recorded By" column by scanning header row 1 (fallback to G/7).
$headerRow = 1
$colIndex = 7
$lastCol = $usedRange.Columns.Count + $usedRange.Column - 1
for ($c = 1; $c -le $lastCol; $c++) {
    $headerVal = $ws.Cells.Item($headerRow, $c).Value()
    if ($headerVal -eq "Recorded By") {
        $colIndex = $c
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversedParts = $parts[($parts.Count - 1)..0]
            $newVal = [string]::Join(", ", $reversedParts)
            if ($val.CompareTo($newVal) -ne 0) {
                $cell.Value = $newVal
            }
        }
    }
}
